# "add feature delete exam" - populate Sheet1 with the exam schedule table
# that the delete-exam feature will operate on (examid, courseid, class,
# date, starttime, endtime, slot), replacing the old placeholder content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "examid"
$ws.Range("B1").Value = "courseid"
$ws.Range("C1").Value = "class"
$ws.Range("D1").NumberFormat = "@"
$ws.Range("D1").Value = "date"
$ws.Range("E1").Value = "starttime"
$ws.Range("F1").Value = "endtime"
$ws.Range("G1").Value = "slot"

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1/12/2000"
$ws.Range("E2").Value = "10:00Am"
$ws.Range("F2").Value = "12:00Am"
$ws.Range("G2").Value = 30

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 2
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1/12/2000"
$ws.Range("E3").Value = "10:00Am"
$ws.Range("F3").Value = "12:00Am"
$ws.Range("G3").Value = 30

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1/1/2000"
$ws.Range("E4").Value = "10:00Am"
$ws.Range("F4").Value = "12:10Am"
$ws.Range("G4").Value = 30

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 2
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1/1/2000"
$ws.Range("E5").Value = "10:00Am"
$ws.Range("F5").Value = "12:00Am"
$ws.Range("G5").Value = 30

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 3
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1/1/2000"
$ws.Range("E6").Value = "10:00Am"
$ws.Range("F6").Value = "12:00Am"
$ws.Range("G6").Value = 30

# Widen the date column to fit its text content.
$ws.Columns.Item(4).ColumnWidth = 9.6

# Printer/page setup used for this sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the selection where the author left it.
[void]$ws.Range("D10").Select()
